# Auto-generated script to apply 2025-05-27 daily crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 2524
$ws.Range("L3").Value = 2539
$ws.Range("E4").Value = 2052
$ws.Range("L4").Value = 690
$ws.Range("L5").Value = 151
$ws.Range("K6").Value = 9123
$ws.Range("L6").Value = 2296
$ws.Range("E7").Value = 26057
$ws.Range("K7").Value = 27559
$ws.Range("L7").Value = 8200

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 265
$ws.Range("L8").Value = 521
$ws.Range("L10").Value = 53
$ws.Range("L11").Value = 144
$ws.Range("L14").Value = 40
$ws.Range("L16").Value = 18
$ws.Range("L17").Value = 14
$ws.Range("L18").Value = 60
$ws.Range("L19").Value = 232
$ws.Range("L20").Value = 211
$ws.Range("L22").Value = 27
$ws.Range("L26").Value = 7
$ws.Range("L27").Value = 86
$ws.Range("L29").Value = 427
$ws.Range("L33").Value = 366
$ws.Range("L37").Value = 299
$ws.Range("L41").Value = 39
$ws.Range("L42").Value = 265
$ws.Range("L44").Value = 60
$ws.Range("L45").Value = 15
$ws.Range("L46").Value = 17
$ws.Range("L47").Value = 62
$ws.Range("L51").Value = 92
$ws.Range("L52").Value = 162
$ws.Range("L53").Value = 101
$ws.Range("L54").Value = 165
$ws.Range("L55").Value = 76
$ws.Range("L60").Value = 48
$ws.Range("E63").Value = 386
$ws.Range("K63").Value = 157
$ws.Range("L63").Value = 27
$ws.Range("L65").Value = 153
$ws.Range("L67").Value = 302
$ws.Range("L71").Value = 22
$ws.Range("L72").Value = 34
$ws.Range("L73").Value = 66
$ws.Range("L75").Value = 34
$ws.Range("L76").Value = 99
$ws.Range("L77").Value = 50
$ws.Range("L78").Value = 106
$ws.Range("L79").Value = 223
$ws.Range("L80").Value = 29
$ws.Range("L84").Value = 83
$ws.Range("L85").Value = 428
$ws.Range("L88").Value = 110
$ws.Range("L89").Value = 101
$ws.Range("L90").Value = 80
$ws.Range("L91").Value = 119
$ws.Range("L94").Value = 100
$ws.Range("L96").Value = 81
$ws.Range("L97").Value = 74
$ws.Range("L99").Value = 133
$ws.Range("L100").Value = 14
$ws.Range("E101").Value = 26057
$ws.Range("K101").Value = 27559
$ws.Range("L101").Value = 8200

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 77
$ws.Range("L3").Value = 85
$ws.Range("L7").Value = 265

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 52
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 35
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 172
$ws.Range("L7").Value = 428

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 52
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 162

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 151
$ws.Range("L3").Value = 173
$ws.Range("L6").Value = 139
$ws.Range("L7").Value = 521

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 113
$ws.Range("L7").Value = 366

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L6").Value = 92
$ws.Range("L7").Value = 299

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 58
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 104
$ws.Range("L5").Value = 8
$ws.Range("L7").Value = 302

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 30
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 83

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 36
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 137
$ws.Range("L3").Value = 155
$ws.Range("L6").Value = 112
$ws.Range("L7").Value = 427

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 75
$ws.Range("L7").Value = 232

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 17
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 81
$ws.Range("L7").Value = 265

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 47
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 80
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 66
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 211

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 100

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("L2").Value = 1
$ws.Range("L7").Value = 7

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 31
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L4").Value = 2
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 18
